$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 20 (C20 and D20 change from 5 to 4)
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4

# Update the selected/active cell to E20 (was E10)
$ws.Range("E20").Select()
